$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-D (Date, Time, Weekday, Week) hold text-like values
# ("2024-01-08", "17:25:10", "Monday", "01"). Excel's automatic type
# inference would otherwise turn these into date/time serials or numbers,
# so force a Text number format before assigning, then clear the
# formatting again (without touching the values) so the new cells end up
# with the same "no explicit style" look as the rest of the sheet.
$textRange = $ws.Range("A34:D34")
$textRange.NumberFormat = "@"

$ws.Range("A34").Value = "2024-01-08"
$ws.Range("B34").Value = "17:25:10"
$ws.Range("C34").Value = "Monday"
$ws.Range("D34").Value = "01"

$textRange.ClearFormats()

$ws.Range("E34").Value = 139524
$ws.Range("F34").Value = 142963
$ws.Range("G34").Value = 172408
$ws.Range("H34").Value = 147268
$ws.Range("I34").Value = -1
$ws.Range("J34").Value = 118186
$ws.Range("K34").Value = 224741
$ws.Range("L34").Value = 249750
$ws.Range("M34").Value = 185117
$ws.Range("N34").Value = 110390
$ws.Range("O34").Value = 40638
$ws.Range("P34").Value = 30802
$ws.Range("Q34").Value = 72433
$ws.Range("R34").Value = -1
$ws.Range("S34").Value = 42121
$ws.Range("T34").Value = -1
